$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Rows 27 and 28 swap coins (Monero <-> ImmutableX) with updated price/volume
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D27") "1.455"
$ws.Range("E27").Value = "  +18.71%  "

$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D28") "151.96"
$ws.Range("E28").Value = "  +2.19%  "

# Price and Volume(1h) updates for the remaining rows
$ws.Range("D2").Value = "26.473.41"
$ws.Range("E2").Value = "  +6.71%  "
$ws.Range("D3").Value = "1.731.39"
$ws.Range("E3").Value = "  +4.40%  "
Set-TextValue $ws.Range("D4") "0.9958"
$ws.Range("E4").Value = "  -0.43%  "
Set-TextValue $ws.Range("D5") "333.88"
$ws.Range("E5").Value = "  +5.40%  "
Set-TextValue $ws.Range("D6") "0.9955"
$ws.Range("E6").Value = "  -0.20%  "
Set-TextValue $ws.Range("D7") "0.3715"
$ws.Range("E7").Value = "  +2.31%  "
Set-TextValue $ws.Range("D8") "49.12"
$ws.Range("E8").Value = "  +4.80%  "
Set-TextValue $ws.Range("D9") "0.3381"
$ws.Range("E9").Value = "  +3.53%  "
Set-TextValue $ws.Range("D10") "1.203"
$ws.Range("E10").Value = "  +5.53%  "
Set-TextValue $ws.Range("D11") "0.07504"
$ws.Range("E11").Value = "  +6.25%  "
Set-TextValue $ws.Range("D12") "0.9970"
$ws.Range("E12").Value = "  -0.04%  "
Set-TextValue $ws.Range("D13") "6.379"
$ws.Range("E13").Value = "  +5.49%  "
Set-TextValue $ws.Range("D14") "20.44"
$ws.Range("E14").Value = "  +4.36%  "
Set-TextValue $ws.Range("D15") "6.982"
$ws.Range("E15").Value = "  +5.41%  "
$ws.Range("D16").Value = "1.713.75"
$ws.Range("E16").Value = "  +3.04%  "
Set-TextValue $ws.Range("D17") "0.00001087"
$ws.Range("E17").Value = "  +3.71%  "
Set-TextValue $ws.Range("D18") "0.06691"
$ws.Range("E18").Value = "  +1.02%  "
Set-TextValue $ws.Range("D19") "82.94"
$ws.Range("E19").Value = "  +4.71%  "
Set-TextValue $ws.Range("D20") "0.9962"
$ws.Range("E20").Value = "  -0.12%  "
Set-TextValue $ws.Range("D21") "16.69"
$ws.Range("E21").Value = "  +5.91%  "
Set-TextValue $ws.Range("D22") "6.157"
$ws.Range("E22").Value = "  +3.88%  "
Set-TextValue $ws.Range("D23") "13.18"
$ws.Range("E23").Value = "  +4.98%  "
$ws.Range("D24").Value = "26.332.88"
$ws.Range("E24").Value = "  +6.27%  "
Set-TextValue $ws.Range("D25") "2.465"
$ws.Range("E25").Value = "  +1.30%  "
Set-TextValue $ws.Range("D26") "2.495"
$ws.Range("E26").Value = "  +3.83%  "
Set-TextValue $ws.Range("D29") "19.48"
$ws.Range("E29").Value = "  +4.68%  "
$ws.Range("D30").Value = "1.905.51"
$ws.Range("E30").Value = "  +3.13%  "
Set-TextValue $ws.Range("D31") "131.18"
$ws.Range("E31").Value = "  +4.32%  "
Set-TextValue $ws.Range("D32") "4.116"
$ws.Range("E32").Value = "  +0.85%  "
Set-TextValue $ws.Range("D33") "6.076"
$ws.Range("E33").Value = "  +4.22%  "
Set-TextValue $ws.Range("D34") "0.08574"
$ws.Range("E34").Value = "  +1.44%  "
Set-TextValue $ws.Range("D35") "1.712"
$ws.Range("E35").Value = "  +2.36%  "
Set-TextValue $ws.Range("D36") "13.16"
$ws.Range("E36").Value = "  +7.20%  "
Set-TextValue $ws.Range("D37") "5.440"
$ws.Range("E37").Value = "  +4.55%  "
Set-TextValue $ws.Range("D38") "0.02352"
$ws.Range("E38").Value = "  +5.08%  "
Set-TextValue $ws.Range("D39") "0.06336"
$ws.Range("E39").Value = "  +4.79%  "
Set-TextValue $ws.Range("D40") "8.694"
$ws.Range("E40").Value = "  +5.68%  "
Set-TextValue $ws.Range("D41") "0.2159"
$ws.Range("E41").Value = "  +4.11%  "
Set-TextValue $ws.Range("D42") "1.239"
$ws.Range("E42").Value = "  -3.11%  "
Set-TextValue $ws.Range("D43") "0.6234"
$ws.Range("E43").Value = "  +5.21%  "
Set-TextValue $ws.Range("D44") "14.37"
$ws.Range("E44").Value = "  +12.10%  "
Set-TextValue $ws.Range("D45") "0.9971"
$ws.Range("E45").Value = "  +0.03%  "
Set-TextValue $ws.Range("D46") "3.869"
$ws.Range("E46").Value = "  +1.24%  "
Set-TextValue $ws.Range("D47") "0.6035"
$ws.Range("E47").Value = "  +6.81%  "
Set-TextValue $ws.Range("D48") "129.46"
$ws.Range("E48").Value = "  +3.51%  "
Set-TextValue $ws.Range("D49") "2.052"
$ws.Range("E49").Value = "  +4.88%  "
Set-TextValue $ws.Range("D50") "0.07342"
$ws.Range("E50").Value = "  +4.99%  "
Set-TextValue $ws.Range("D51") "77.58"
$ws.Range("E51").Value = "  +3.42%  "
